$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a new value that looks numeric must be forced to Text format first,
# otherwise Excel auto-converts strings like "519.90" or "1.00" into numbers and
# drops significant trailing zeros / reformats them (e.g. "1.00" -> 1, "0.160" -> 0.16).
$textFormatCells = @("D5","D6","D8","D9","D15","D16","D19","D20","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D35","D37","D38","D39","D40","D41","D42","D44","D45","D46","D48","D50")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row, in the same order as the source diff.
# Row 2
$ws.Range("D2").Value = '58.129.50'

# Row 3
$ws.Range("D3").Value = '2.474.14'
$ws.Range("E3").Value = '  -0.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = '519.90'
$ws.Range("E5").Value = '  -2.92%  '

# Row 6
$ws.Range("D6").Value = '131.25'
$ws.Range("E6").Value = '  -3.73%  '

# Row 7
$ws.Range("E7").Value = '  +0.09%  '

# Row 8
$ws.Range("D8").Value = '0.559'
$ws.Range("E8").Value = '  -1.27%  '

# Row 9
$ws.Range("D9").Value = '0.0994'
$ws.Range("E9").Value = '  -1.42%  '

# Row 10
$ws.Range("E10").Value = '  -0.43%  '

# Row 11
$ws.Range("E11").Value = '  +0.28%  '

# Row 12
$ws.Range("E12").Value = '  -0.87%  '

# Row 13
$ws.Range("D13").Value = '2.914.54'
$ws.Range("E13").Value = '  -1.61%  '

# Row 14
$ws.Range("D14").Value = '58.073.42'
$ws.Range("E14").Value = '  -1.49%  '

# Row 15
$ws.Range("D15").Value = '22.04'
$ws.Range("E15").Value = '  -3.83%  '

# Row 16
$ws.Range("D16").Value = '0.0000136'
$ws.Range("E16").Value = '  -1.69%  '

# Row 17
$ws.Range("D17").Value = '2.481.20'
$ws.Range("E17").Value = '  -1.38%  '

# Row 18
$ws.Range("E18").Value = '  -2.28%  '

# Row 19
$ws.Range("D19").Value = '4.17'
$ws.Range("E19").Value = '  -2.38%  '

# Row 20
$ws.Range("D20").Value = '319.14'
$ws.Range("E20").Value = '  -1.29%  '

# Row 21
$ws.Range("E21").Value = '  -0.06%  '

# Row 22
$ws.Range("D22").Value = '5.74'
$ws.Range("E22").Value = '  -2.87%  '

# Row 23
$ws.Range("D23").Value = '64.12'
$ws.Range("E23").Value = '  -1.77%  '

# Row 24
$ws.Range("D24").Value = '0.409'
$ws.Range("E24").Value = '  -2.53%  '

# Row 25
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.08%  '

# Row 26
$ws.Range("D26").Value = '0.160'
$ws.Range("E26").Value = '  -3.28%  '

# Row 27
$ws.Range("D27").Value = '7.37'
$ws.Range("E27").Value = '  -2.78%  '

# Row 29
$ws.Range("D29").Value = '1.70'
$ws.Range("E29").Value = '  -3.62%  '

# Row 30
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '166.65'
$ws.Range("E30").Value = '  +0.11%  '

# Row 31
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '6.31'
$ws.Range("E31").Value = '  -5.78%  '

# Row 32
$ws.Range("D32").Value = '1.17'
$ws.Range("E32").Value = '  +0.84%  '

# Row 33
$ws.Range("E33").Value = '  -0.08%  '

# Row 34
$ws.Range("E34").Value = '  +0.09%  '

# Row 35
$ws.Range("D35").Value = '18.07'
$ws.Range("E35").Value = '  -1.73%  '

# Row 36
$ws.Range("E36").Value = '  -10.37%  '

# Row 37
$ws.Range("D37").Value = '3.97'
$ws.Range("E37").Value = '  -2.69%  '

# Row 38
$ws.Range("D38").Value = '1.48'
$ws.Range("E38").Value = '  -3.58%  '

# Row 39
$ws.Range("D39").Value = '0.790'
$ws.Range("E39").Value = '  -2.57%  '

# Row 40
$ws.Range("D40").Value = '3.47'
$ws.Range("E40").Value = '  -3.77%  '

# Row 41
$ws.Range("D41").Value = '275.65'
$ws.Range("E41").Value = '  -3.46%  '

# Row 42
$ws.Range("D42").Value = '5.02'
$ws.Range("E42").Value = '  -4.01%  '

# Row 43
$ws.Range("E43").Value = '  -1.40%  '

# Row 44
$ws.Range("D44").Value = '126.44'
$ws.Range("E44").Value = '  -4.66%  '

# Row 45
$ws.Range("D45").Value = '0.0906'
$ws.Range("E45").Value = '  -2.01%  '

# Row 46
$ws.Range("D46").Value = '0.0489'
$ws.Range("E46").Value = '  -3.50%  '

# Row 47
$ws.Range("E47").Value = '  -2.36%  '

# Row 48
$ws.Range("D48").Value = '17.15'
$ws.Range("E48").Value = '  -0.35%  '

# Row 49
$ws.Range("D49").Value = '1.738.19'
$ws.Range("E49").Value = '  -1.84%  '

# Row 50
$ws.Range("D50").Value = '0.973'
$ws.Range("E50").Value = '  -1.98%  '

# Row 51
$ws.Range("E51").Value = '  -0.97%  '

